$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Updated visitor-home-location agency figures (recomputed with new
# agency boundary polygons built via Tigris).
$ws.Range("B8").Value = 31759
$ws.Range("C8").Value = 5.03811416273365

$ws.Range("B9").Value = 8716
$ws.Range("C9").Value = 1.38276494775071

$ws.Range("B10").Value = 9522
$ws.Range("C10").Value = 1.51062703476646

$ws.Range("B11").Value = 9285
$ws.Range("C11").Value = 1.47292319366238

$ws.Range("B12").Value = 229892
$ws.Range("C12").Value = 36.4696515855179

$ws.Range("B13").Value = 21216
$ws.Range("C13").Value = 3.36567549861431

$ws.Range("C14").Value = 1.06321126436847

$ws.Range("B15").Value = 27535
$ws.Range("C15").Value = 4.36809854758978

$ws.Range("B16").Value = 252567
$ws.Range("C16").Value = 40.0667248278992

$ws.Range("B17").Value = 5979
$ws.Range("C17").Value = 0.94846521487083
